# Update the "6) Please clarify/correct" reply note (C5) with a concise
# summary of the forward-progress explanation, and restyle that cell so it
# reads in the normal body font (Calibri (Body), black) and is centered
# both horizontally and vertically, matching the rest of the "Corrections"
# column replies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C5")
$cell.Value = "Add Section 3.2 and Figure 3.3"

$cell.Font.Name = "Calibri (Body)"
$cell.Font.Size = 12
$cell.Font.ThemeColor = 1
$cell.HorizontalAlignment = -4108  # xlCenter
$cell.VerticalAlignment = -4108    # xlCenter
